# Auto-generated: apply market-price data refresh to Sephirot_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1680.6552
$ws.Range("I15").Value = 1680.6552
$ws.Range("K15").Value = 5041.9656
$ws.Range("M15").Value = -4872.9656
$ws.Range("H41").Value = 495.25
$ws.Range("I41").Value = 493.66666
$ws.Range("J41").Value = 500
$ws.Range("K41").Value = 493.66666
$ws.Range("L41").Value = 500
$ws.Range("M41").Value = -53.66665999999998
$ws.Range("N41").Value = -1380
$ws.Range("H80").Value = 2700
$ws.Range("I80").Value = 2700
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 8100
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -7102
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 2700
$ws.Range("I83").Value = 2700
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 24300
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -19308
$ws.Range("N83").ClearContents()
$ws.Range("H92").Value = 967.94116
$ws.Range("I92").Value = 1002.8571
$ws.Range("J92").Value = 805
$ws.Range("K92").Value = 1002.8571
$ws.Range("L92").Value = 805
$ws.Range("M92").Value = 245.1429000000001
$ws.Range("N92").Value = -3301
$ws.Range("H137").Value = 2569.8572
$ws.Range("I137").Value = 1747.5
$ws.Range("K137").Value = 5242.5
$ws.Range("M137").Value = -2692.5
$ws.Range("H138").Value = 3828
$ws.Range("I138").Value = 3249.75
$ws.Range("J138").Value = 3910.6072
$ws.Range("K138").Value = 9749.25
$ws.Range("L138").Value = 11731.8216
$ws.Range("M138").Value = -4609.25
$ws.Range("N138").Value = -22011.8216

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3545.1516
$ws.Range("I32").Value = 3511.9
$ws.Range("J32").Value = 3877.6667
$ws.Range("K32").Value = 3511.9
$ws.Range("L32").Value = 3877.6667
$ws.Range("M32").Value = -3224.9
$ws.Range("N32").Value = -4451.6667
$ws.Range("H74").Value = 1697.3334
$ws.Range("I74").Value = 1697.3334
$ws.Range("K74").Value = 1697.3334
$ws.Range("M74").Value = -823.3334
$ws.Range("H77").Value = 1697.3334
$ws.Range("I77").Value = 1697.3334
$ws.Range("K77").Value = 8486.666999999999
$ws.Range("M77").Value = -4118.666999999999
$ws.Range("H122").Value = 6845.3887
$ws.Range("I122").Value = 6914.467
$ws.Range("K122").Value = 20743.401
$ws.Range("M122").Value = -18293.401
$ws.Range("H132").Value = 1866.1212
$ws.Range("I132").Value = 1022.46155
$ws.Range("K132").Value = 3067.38465
$ws.Range("M132").Value = -537.38465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4401.143
$ws.Range("I20").Value = 4159.8
$ws.Range("K20").Value = 4159.8
$ws.Range("M20").Value = -3912.8
$ws.Range("H86").Value = 1860.4286
$ws.Range("I86").Value = 1860.4286
$ws.Range("K86").Value = 1860.4286
$ws.Range("M86").Value = -737.4286
$ws.Range("H89").Value = 1860.4286
$ws.Range("I89").Value = 1860.4286
$ws.Range("K89").Value = 9302.143
$ws.Range("M89").Value = -3686.143
$ws.Range("H106").Value = 242333.33
$ws.Range("J106").Value = 242333.33
$ws.Range("L106").Value = 242333.33
$ws.Range("N106").Value = -244857.33

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 10000000
$ws.Range("I3").Value = 10000000
$ws.Range("K3").Value = 10000000
$ws.Range("M3").Value = -9999887
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H62").Value = 3873.625
$ws.Range("I62").Value = 3496.6
$ws.Range("K62").Value = 3496.6
$ws.Range("M62").Value = -2872.6
$ws.Range("H65").Value = 3873.625
$ws.Range("I65").Value = 3496.6
$ws.Range("K65").Value = 17483
$ws.Range("M65").Value = -14363
$ws.Range("H132").Value = 2381.4666
$ws.Range("I132").Value = 2194.8
$ws.Range("J132").Value = 2754.8
$ws.Range("K132").Value = 6584.400000000001
$ws.Range("L132").Value = 8264.400000000001
$ws.Range("M132").Value = -4054.400000000001
$ws.Range("N132").Value = -13324.4
$ws.Range("H134").Value = 1738.1111
$ws.Range("I134").Value = 1642.875
$ws.Range("K134").Value = 4928.625
$ws.Range("M134").Value = -2393.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 73.333336
$ws.Range("I12").Value = 100
$ws.Range("K12").Value = 300
$ws.Range("M12").Value = -127
$ws.Range("H14").Value = 1424.5
$ws.Range("I14").Value = 1424.5
$ws.Range("K14").Value = 4273.5
$ws.Range("M14").Value = -4100.5
$ws.Range("H107").Value = 1013.1429
$ws.Range("J107").Value = 1078.6
$ws.Range("L107").Value = 3235.8
$ws.Range("N107").Value = -7075.799999999999
$ws.Range("H122").Value = 2961.6365
$ws.Range("I122").Value = 655.5714
$ws.Range("J122").Value = 6997.25
$ws.Range("K122").Value = 5900.1426
$ws.Range("L122").Value = 62975.25
$ws.Range("M122").Value = -3450.1426
$ws.Range("N122").Value = -67875.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10000
$ws.Range("I70").Value = 10000
$ws.Range("K70").Value = 10000
$ws.Range("M70").Value = -9730
$ws.Range("H73").Value = 10000
$ws.Range("I73").Value = 10000
$ws.Range("K73").Value = 10000
$ws.Range("M73").Value = -9064
$ws.Range("H102").Value = 2691.5715
$ws.Range("I102").Value = 2307.6365
$ws.Range("K102").Value = 2307.6365
$ws.Range("M102").Value = -685.6365000000001
$ws.Range("H122").Value = 3800.75
$ws.Range("I122").Value = 4067.6667
$ws.Range("K122").Value = 12203.0001
$ws.Range("M122").Value = -9753.000100000001
$ws.Range("H126").Value = 3200
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3800
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 11400
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -16340
$ws.Range("H134").Value = 51726668
$ws.Range("J134").Value = 51726668
$ws.Range("L134").Value = 155180004
$ws.Range("N134").Value = -155185074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3949.5
$ws.Range("I7").Value = 3949.5
$ws.Range("K7").Value = 3949.5
$ws.Range("M7").Value = -3837.5
$ws.Range("H22").Value = 5466.3335
$ws.Range("I22").Value = 1874.5
$ws.Range("J22").Value = 6772.4546
$ws.Range("K22").Value = 1874.5
$ws.Range("L22").Value = 6772.4546
$ws.Range("M22").Value = -1579.5
$ws.Range("N22").Value = -7362.4546
$ws.Range("H27").Value = 5466.3335
$ws.Range("I27").Value = 1874.5
$ws.Range("J27").Value = 6772.4546
$ws.Range("K27").Value = 1874.5
$ws.Range("L27").Value = 6772.4546
$ws.Range("M27").Value = -1767.5
$ws.Range("N27").Value = -6986.4546
$ws.Range("H40").Value = 3000
$ws.Range("I40").Value = 3000
$ws.Range("K40").Value = 3000
$ws.Range("M40").Value = -2864
$ws.Range("H100").Value = 2750
$ws.Range("I100").Value = 2750
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2750
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2209
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 3949.5
$ws.Range("I126").Value = 3949.5
$ws.Range("K126").Value = 11848.5
$ws.Range("M126").Value = -9378.5
$ws.Range("H136").Value = 1397
$ws.Range("I136").Value = 1397
$ws.Range("K136").Value = 4191
$ws.Range("M136").Value = -1641

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 400
$ws.Range("I100").Value = 400
$ws.Range("K100").Value = 800
$ws.Range("M100").Value = -259
$ws.Range("H122").Value = 1014.2222
$ws.Range("I122").Value = 1174.8334
$ws.Range("J122").Value = 693
$ws.Range("K122").Value = 3524.5002
$ws.Range("L122").Value = 2079
$ws.Range("M122").Value = -1074.5002
$ws.Range("N122").Value = -6979
$ws.Range("H126").Value = 1570.7142
$ws.Range("I126").Value = 998.75
$ws.Range("K126").Value = 2996.25
$ws.Range("M126").Value = -526.25
$ws.Range("H130").Value = 47500
$ws.Range("J130").Value = 47500
$ws.Range("L130").Value = 47500
$ws.Range("N130").Value = -57540
$ws.Range("H132").Value = 4630.6665
$ws.Range("I132").Value = 2972.6667
$ws.Range("J132").Value = 5874.1665
$ws.Range("K132").Value = 8918.000100000001
$ws.Range("L132").Value = 17622.4995
$ws.Range("M132").Value = -6388.000100000001
$ws.Range("N132").Value = -22682.4995
$ws.Range("H136").Value = 2024.375
$ws.Range("I136").Value = 2024.375
$ws.Range("K136").Value = 6073.125
$ws.Range("M136").Value = -3523.125
